$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (surgical substring edits, preserve surrounding text) ---
$ws.Range("A8").Characters(21,1).Text = "8"
$ws.Range("C9").Characters(27,9).Text = "2/17/2025"
$ws.Range("C9").Characters(47,9).Text = "2/23/2025"

# --- Data table updates (rows 14-31) ---
$ws.Range("C14").Value = 1
$ws.Range("J14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("G14").PasteSpecial(-4122)
$ws.Range("H14").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("H14").PasteSpecial(-4122)
$ws.Range("I14").Value = 2
$ws.Range("K14").Value = 100
$ws.Range("L14").Value = 100
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = -50
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 4
$ws.Range("H15").Value = -25
$ws.Range("I15").Value = 4
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = -20
$ws.Range("L15").Value = 33.333333333333
$ws.Range("M15").Value = 100
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -25
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = 11.111111111111
$ws.Range("I16").Value = 41
$ws.Range("J16").Value = 44
$ws.Range("K16").Value = -6.818181818181
$ws.Range("L16").Value = 10.810810810810
$ws.Range("M16").Value = -12.765957446808
$ws.Range("N16").Value = -68.217054263565
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 20
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = -3.448275862068
$ws.Range("I17").Value = 56
$ws.Range("J17").Value = 57
$ws.Range("K17").Value = -1.754385964912
$ws.Range("L17").Value = 16.666666666666
$ws.Range("M17").Value = 60
$ws.Range("N17").Value = 154.545454545455
$ws.Range("C18").Value = 6
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 36
$ws.Range("G18").Value = 13
$ws.Range("H18").Value = 176.923076923077
$ws.Range("I18").Value = 54
$ws.Range("J18").Value = 25
$ws.Range("K18").Value = 116
$ws.Range("L18").Value = 63.636363636363
$ws.Range("M18").Value = 14.893617021276
$ws.Range("N18").Value = -79.069767441860
$ws.Range("C19").Value = 24
$ws.Range("D19").Value = 28
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 64
$ws.Range("G19").Value = 83
$ws.Range("H19").Value = -22.891566265060
$ws.Range("I19").Value = 119
$ws.Range("J19").Value = 149
$ws.Range("K19").Value = -20.134228187919
$ws.Range("L19").Value = 48.75
$ws.Range("M19").Value = 164.444444444444
$ws.Range("N19").Value = 50.632911392405
$ws.Range("C20").Value = 8
$ws.Range("E20").Value = 14.285714285714
$ws.Range("G20").Value = 37
$ws.Range("H20").Value = -2.702702702702
$ws.Range("I20").Value = 74
$ws.Range("J20").Value = 70
$ws.Range("K20").Value = 5.714285714285
$ws.Range("L20").Value = 7.246376811594
$ws.Range("M20").Value = 138.709677419355
$ws.Range("N20").Value = -71.206225680933
$ws.Range("C21").Value = 49
$ws.Range("D21").Value = 48
$ws.Range("E21").Value = 2.083333333333
$ws.Range("F21").Value = 189
$ws.Range("G21").Value = 184
$ws.Range("H21").Value = 2.717391304347
$ws.Range("I21").Value = 350
$ws.Range("J21").Value = 351
$ws.Range("K21").Value = -0.284900284900
$ws.Range("L21").Value = 29.151291512915
$ws.Range("M21").Value = 68.269230769230
$ws.Range("N21").Value = -53.519256308100
$ws.Range("C22").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$ws.Range("J14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("E15").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 4
$ws.Range("K22").Value = 0
$ws.Range("C23").Value = 2
$ws.Range("J14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 9
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 18
$ws.Range("J23").Value = 20
$ws.Range("K23").Value = -10
$ws.Range("L23").Value = -5.263157894736
$ws.Range("M23").Value = 80
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = -18.181818181818
$ws.Range("F24").Value = 91
$ws.Range("G24").Value = 139
$ws.Range("H24").Value = -34.532374100719
$ws.Range("I24").Value = 182
$ws.Range("J24").Value = 235
$ws.Range("K24").Value = -22.553191489361
$ws.Range("L24").Value = -29.457364341085
$ws.Range("M24").Value = 35.820895522388
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -55.555555555555
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 69
$ws.Range("H25").Value = -62.318840579710
$ws.Range("I25").Value = 66
$ws.Range("J25").Value = 115
$ws.Range("K25").Value = -42.608695652173
$ws.Range("L25").Value = -43.103448275862
$ws.Range("C26").Value = 13
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 62.5
$ws.Range("F26").Value = 62
$ws.Range("G26").Value = 43
$ws.Range("H26").Value = 44.186046511627
$ws.Range("I26").Value = 99
$ws.Range("J26").Value = 73
$ws.Range("K26").Value = 35.616438356164
$ws.Range("L26").Value = 30.263157894736
$ws.Range("M26").Value = 25.316455696202
$ws.Range("C27").Value = 2
$ws.Range("E27").Value = 100
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 5
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = -16.666666666666
$ws.Range("L27").Value = -28.571428571428
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 2
$ws.Range("J14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = -50
$ws.Range("E15").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 80
$ws.Range("I28").Value = 14
$ws.Range("J28").Value = 9
$ws.Range("K28").Value = 55.555555555555
$ws.Range("L28").Value = 100
$ws.Range("C29").Value = 2
$ws.Range("J14").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("D29").Value = 1
$ws.Range("J14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("E29").Value = 100
$ws.Range("E15").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 1
$ws.Range("J14").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("H29").Value = 200
$ws.Range("E15").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("I29").Value = 3
$ws.Range("J29").Value = 2
$ws.Range("K29").Value = 50
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -40
$ws.Range("C30").Value = 1
$ws.Range("J14").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("D30").Value = 1
$ws.Range("J14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = 0
$ws.Range("E15").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("F30").Value = 2
$ws.Range("G30").Value = 1
$ws.Range("J14").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("H30").Value = 100
$ws.Range("E15").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("I30").Value = 2
$ws.Range("J30").Value = 2
$ws.Range("L30").Value = -33.333333333333
$ws.Range("M30").Value = -33.333333333333
$ws.Range("N30").Value = -50
$ws.Range("C31").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("C31").PasteSpecial(-4122)
